$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "ng-Paths" result columns (J/K/L)
$ws.Range("J1").Value = "ng-Paths"
$ws.Range("K1").Value = 8

# Fix typo in instance name
$ws.Range("A3").Value = "E-n23-k3"

$ws.Range("K3").Value = "558,6…"
$ws.Range("K2").Value = "373,5…"
$ws.Range("L3").Value = "190s"

$ws.Range("K4").Value = "481,0…"
$ws.Range("L4").Value = "292s"

$ws.Range("K7").Value = "1002,2…"
$ws.Range("L7").Value = "324s"

$ws.Range("K8").Select()
